$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-missing hours logged for the 2017-08-15 entry
# (row 22, column D) -> this also causes F3's SUM(D3:D33) to recalc from
# 78.5 to 86.5.
$ws.Range("D22").Value = 8

# Scroll the view down a bit (so row 4 is at the top of the visible
# window) and move the active selection/cursor to F20.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("F20").Select()

$wb.Save()
